# Apply all edits described by the diff using Find (match-only) + Range.Text assignment.
# Using Range.Text (instead of the Find.Execute "Replace" parameter) avoids Word's
# "smart quotes" AutoCorrect/AutoFormat from mangling straight apostrophes in the new text.
$d = $word.ActiveDocument

# Replacement 1
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Ativação: 01/01/2023", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 1 found:" $found1
if ($found1) {
    $rng1.Text = "Ativação: 01/01/2024"
}

# Replacement 2
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("1341653 - Maria José Ramos Sandim^l", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 2 found:" $found2
if ($found2) {
    $rng2.Text = ""
}

# Replacement 3
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Óptica de raios; Ondas eletromagnéticas: fase e polarização; Interferência; Coerência; Difração; Óptica de Fourier; Interação da luz com a matéria; Guias de ondas metálicos e dielétricos; Óptica de cristais; Óptica não linear.", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 3 found:" $found3
if ($found3) {
    $rng3.Text = "Descrição ondulatória e quântica da luz. Propriedades da luz. Interação da luz com a matéria. Aplicações."
}

# Replacement 4
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Ray optics; Electromagnetic waves: phase and polarization; Interference; Coherence; Diffraction; Fourier optics; Interaction of light with matter; Metallic and dielectric waveguides; Crystal optics; Non-linear optics.", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 4 found:" $found4
if ($found4) {
    $rng4.Text = "Presentation of the wave and quantum description of light, study of the properties of light, the interaction of light with matter and applications of physical optics."
}

# Replacement 5
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("Óptica de raios. Introdução. Propagação de luz em meios homogêneos. Propagação de luz em meios não homogêneos. A lei de Snell generalizada. O princípio de Fermat. A equação dos raios. A função eikonal. Analogia ente a mecânica clássica e a óptica geométrica. O potencial óptico.Ondas eletromagnéticas. Ondas harmônicas unidimensionais. Ondas planas e esféricas. Ondas gaussianas. Propagação do feixe gaussiano. Vetor de Poynting. Intensidade.A fase da onda eletromagnética. Velocidades de fase e de grupo. Dispersão. Efeito Doppler. Aplicações astronômicas. Alargamento de linhas espectrais. Óptica relativística. Modulação eletroóptica de frequência. Automodulação de fase. Polarização das ondas eletromagnéticas. Polarização linear. Polarização elíptica. Polarização circular. Obtenção de luz linearmente polarizada. Equações de Fresnel. Polarização por reflexão total interna. Matrizes de Jones. Atividade óptica. Efeito Faraday. Isoladores ópticos. Efeito Pockels. Efeitos Kerr e Cotton-Mouton. Chaveamento eletroóptico.Interferência. Princípio da superposição. Interferência por divisão da frente de onda. Interferência por divisão de amplitudes. Interferômetro de Fabry-Perot. Analisador de espectro óptico. Teoria de películas.Coerência. Introdução. Coerência temporal. Resolução espectral de um trem de ondas finito. Coerência espacial. Medidas de diâmetros de estrelas.Difração. Princípio de Huygens. Fórmula de Fresnel-Kirchhoff. Princípio de Babinet. Difração de Fraunhofer. Difração por uma abertura circular. Rede de difração. Padrões de difração de Fresnel. Óptica de Fourier.  Microscopia por contraste de fase.  Holografia. Interação da radiação com a matéria.  Modelo do oscilador harmônico.  Dispersão cromática do índice de refração. Absorção. Espalhamento Rayleigh. Força da radiação em átomo neutro.Óptica não linear. Susceptibilidade não linear, processos paramétricos e não paramétricos. Geração de freqüências. Casamento de fase.", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 5 found:" $found5
if ($found5) {
    $rng5.Text = "O que é luz? Reflexão. Refração. Difração. Polarização. Formação de imagens: Transformada de Fourier. Ondas eletromagnéticas. Equações de Maxwell. Propagação da luz em diferentes meios: vácuo, dielétrico, condutor. Transporte de energia. Condições de contorno entre diferentes meios: vácuo, dielétrico, condutor. Propagação da luz entre diferentes meios: incidência normal e oblíqua na interface entre meios. Coeficientes de Fresnel. Aplicações da Óptica: holografia, laser, fibras ópticas, materiais eletrocrômicos, metamateriais."
}

# Replacement 6
$rng6 = $d.Content
$found6 = $rng6.Find.Execute("Ray optics. Introduction. Propagation of light in homogeneous media. Propagation of light in non-homogeneous media. Generalized Snell's law. Fermat's principle. The equation of rays. The eikonal function. Analogy between classical mechanics and optics geometric The optical potential.Electromagnetic waves. One-dimensional harmonic waves. Flat and spherical waves. Gaussian waves. Gaussian beam propagation. Poynting vector. Intensity.The phase of the electromagnetic wave. Phase and group speeds. Dispersal. Doppler effect. Astronomical applications. Broadening of spectral lines. Relativistic optics. Electro-optical frequency modulation. Phase automodulation.Polarization of electromagnetic waves. Linear polarization. Elliptical Polarization. Circular polarization. Obtaining linearly polarized light. Fresnel equations. Polarization by total internal reflection. Jones matrices. Optical activity. Faraday effect. Optical isolators. Pockels Effect. Kerr and Cotton-Mouton effects. Electro-optical switching.Interference. Superposition principle. Interference by division of the wavefront. Amplitude division interference. Fabry-Perot interferometer. Optical spectrum analyzer. Film theory.Coherence. Introduction. Temporal coherence. Spectral resolution of a finite wave train. Spatial coherence. Star diameter measurements.Diffraction. Huygens Principle. Fresnel-Kirchhoff formula. Babinet's Principle. Fraunhofer Diffraction. Diffraction through a circular aperture. Diffraction grating. Fresnel diffraction patterns. Fourier optics. Phase contrast microscopy. Holography.Interaction of radiation with matter. Harmonic oscillator model. Chromatic dispersion of the refractive index. Absorption. Rayleigh scattering. Force of radiation on a neutral atom.Non-linear optics. Nonlinear susceptibility, parametric and nonparametric processes. Frequency generation. Phase marriage.", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 6 found:" $found6
if ($found6) {
    $rng6.Text = "What is light? Reflection. Refraction. Diffraction. Polarization. Image formation: Fourier transform. Electromagnetic waves. Maxwell's equations. Propagation of light in different media: vacuum, dielectric, conductor. Energy transport. Boundary conditions between different media: vacuum, dielectric, conductor. Light propagation between different media: normal and oblique incidence at the interface between media. Fresnel coefficients. Optics Applications: holography, laser, optical fibers, electrochromic materials, metamaterials."
}

# Replacement 7
$rng7 = $d.Content
$found7 = $rng7.Find.Execute("HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005.", $true, $false, $false, $false, $false, $true, 1, $false)
Write-Host "Replacement 7 found:" $found7
if ($found7) {
    $rng7.Text = "HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005.J. R. Reitz, F. J. Milford, R. W. Christy, Fundamentos da Teoria Eletromagnética. Editora Campus. 1982."
}
